$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D2:E51 are treated as plain text so numeric-looking strings
# (e.g. "1.00", "26.909.53") are preserved exactly as authored.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "26.909.53"
$ws.Range("E2").Value = "  +0.40%  "

# Row 3
$ws.Range("D3").Value = "1.551.81"
$ws.Range("E3").Value = "  +0.20%  "

# Row 4
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.39%  "

# Row 5
$ws.Range("D5").Value = "206.31"
$ws.Range("E5").Value = "  +0.82%  "

# Row 6
$ws.Range("E6").Value = "  +0.35%  "

# Row 7
$ws.Range("E7").Value = "  +0.38%  "

# Row 8
$ws.Range("B8").Value = "Solana"
$ws.Range("C8").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D8").Value = "21.52"
$ws.Range("E8").Value = "  +0.46%  "

# Row 9
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").Value = "0.247"
$ws.Range("E9").Value = "  +0.68%  "

# Row 10
$ws.Range("E10").Value = "  +0.26%  "

# Row 11
$ws.Range("D11").Value = "0.0858"
$ws.Range("E11").Value = "  +0.11%  "

# Row 12
$ws.Range("D12").Value = "1.772.86"
$ws.Range("E12").Value = "  +0.23%  "

# Row 13
$ws.Range("D13").Value = "1.550.51"
$ws.Range("E13").Value = "  +0.16%  "

# Row 14
$ws.Range("E14").Value = "  +0.60%  "

# Row 15
$ws.Range("E15").Value = "  +0.65%  "

# Row 16
$ws.Range("D16").Value = "26.907.63"
$ws.Range("E16").Value = "  +0.43%  "

# Row 17
$ws.Range("D17").Value = "61.69"
$ws.Range("E17").Value = "  +1.30%  "

# Row 18
$ws.Range("D18").Value = "213.61"
$ws.Range("E18").Value = "  -0.16%  "

# Row 19
$ws.Range("E19").Value = "  +0.38%  "

# Row 20
$ws.Range("D20").Value = "7.24"
$ws.Range("E20").Value = "  -0.16%  "

# Row 21
$ws.Range("D21").Value = "1.01"
$ws.Range("E21").Value = "  +0.41%  "

# Row 22
$ws.Range("E22").Value = "  -0.90%  "

# Row 23
$ws.Range("D23").Value = "9.17"
$ws.Range("E23").Value = "  +1.45%  "

# Row 24
$ws.Range("E24").Value = "  -1.19%  "

# Row 25
$ws.Range("D25").Value = "153.01"
$ws.Range("E25").Value = "  +0.20%  "

# Row 26
$ws.Range("E26").Value = "  +2.33%  "

# Row 27
$ws.Range("D27").Value = "14.86"
$ws.Range("E27").Value = "  -0.25%  "

# Row 28
$ws.Range("E28").Value = "  +0.38%  "

# Row 29
$ws.Range("E29").Value = "  +1.24%  "

# Row 30
$ws.Range("E30").Value = "  -0.55%  "

# Row 31
$ws.Range("E31").Value = "  -0.57%  "

# Row 32
$ws.Range("D32").Value = "3.23"
$ws.Range("E32").Value = "  +1.82%  "

# Row 33
$ws.Range("D33").Value = "1.374.79"
$ws.Range("E33").Value = "  +1.52%  "

# Row 34
$ws.Range("E34").Value = "  +2.04%  "

# Row 35
$ws.Range("E35").Value = "  +3.14%  "

# Row 36
$ws.Range("D36").Value = "0.972"
$ws.Range("E36").Value = "  +6.15%  "

# Row 37
$ws.Range("E37").Value = "  +0.48%  "

# Row 38
$ws.Range("E38").Value = "  +1.08%  "

# Row 39
$ws.Range("D39").Value = "0.523"
$ws.Range("E39").Value = "  +0.00%  "

# Row 40
$ws.Range("E40").Value = "  +0.77%  "

# Row 41
$ws.Range("D41").Value = "1.01"
$ws.Range("E41").Value = "  +0.38%  "

# Row 42
$ws.Range("D42").Value = "0.988"
$ws.Range("E42").Value = "  -0.15%  "

# Row 43
$ws.Range("E43").Value = "  -0.75%  "

# Row 44
$ws.Range("E44").Value = "  +3.22%  "

# Row 45
$ws.Range("D45").Value = "63.63"
$ws.Range("E45").Value = "  +1.25%  "

# Row 46
$ws.Range("E46").Value = "  -2.10%  "

# Row 47
$ws.Range("D47").Value = "1.686.46"
$ws.Range("E47").Value = "  +0.16%  "

# Row 48
$ws.Range("D48").Value = "86.21"
$ws.Range("E48").Value = "  +0.36%  "

# Row 49
$ws.Range("E49").Value = "  -0.11%  "

# Row 50
$ws.Range("D50").Value = "0.0954"
$ws.Range("E50").Value = "  +1.01%  "

# Row 51
$ws.Range("E51").Value = "  +0.63%  "
